# Update the "keprof" template: header row now describes a student's
# NIM / Name / Keprofesian (instead of lecturer's name / kelompok / peminatan).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - field keys
$ws.Range("A1").Value = "nim"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "keprofesian"

# Row 2 - human readable labels
$ws.Range("A2").Value = "NIM Mahasiswa"
$ws.Range("B2").Value = "Nama Mahasiswa"
$ws.Range("C2").Value = "Keprofesian"

# Move/leave the active selection on C2, matching the saved view state.
$ws.Range("C2").Select()
